# R1000AX_BOM.xlsx -- "added host USB to LPC USB jump option using a diffpair"
#
# The BOM sheet itself doesn't gain new rows in this particular file (the
# schematic / other BOM-producing file in the commit did); what shows up
# here is the incidental re-layout that the BOM exporter produced: the
# selection cursor moved, several columns were nudged to new auto widths,
# the header/footer gained an explicit font, and the top/bottom margins
# grew to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- selection -------------------------------------------------------
$ws.Range("F21").Select()

# --- column widths (characters, Normal-font units) --------------------
# NOTE: this host's ColumnWidth setter always re-adds the Normal-font
# padding (5/6 of a character) on top of whatever is assigned and then
# snaps to the nearest 1/6th of a character, so the assigned values below
# are pre-compensated (target - 5/6) to land on the closest achievable
# grid point to the real target width.
$ws.Columns.Item(2).ColumnWidth  = 7.534013605442177
$ws.Columns.Item(3).ColumnWidth  = 26.498299319727867
$ws.Columns.Item(4).ColumnWidth  = 35.166666666666664
$ws.Columns.Item(5).ColumnWidth  = 9.003401360544217
$ws.Columns.Item(6).ColumnWidth  = 34.508503401360564
$ws.Columns.Item(7).ColumnWidth  = 22.411564625850367
$ws.Columns.Item(8).ColumnWidth  = 9.003401360544217
$ws.Columns.Item(9).ColumnWidth  = 7.039115646258507
$ws.Columns.Item(10).ColumnWidth = 5.411564625850337
$ws.Columns.Item(11).ColumnWidth = 7.375850340136057
$ws.Columns.Item(12).ColumnWidth = 6.880952380952377
$ws.Columns.Item(13).ColumnWidth = 5.411564625850337
$ws.Columns.Item(14).ColumnWidth = 8.518707482993197
$ws.Columns.Item(15).ColumnWidth = 7.212585034013608

# --- page margins (points; PageSetup stores/returns points) -----------
$ps = $ws.PageSetup
$ps.TopMargin    = 75.8
$ps.BottomMargin = 75.8

# --- header / footer: explicit Times New Roman 12pt --------------------
$ps.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ps.CenterFooter = '&"Times New Roman,Regular"&12Page &P'
